$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: update the date in A1 (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the price in D29 (960 -> 1570)
$ws.Range("D29").Value = 1570
